$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New species rows appended to the atlas list (stashing species with no eBird data).
# Entry order matches the author's original data-entry order (alphabetical by
# species name: Anous, Aplonis, Charadrius, Pachyptila -> sheet rows 4,3,5,6)
# so the shared-string table is built up in the same sequence.
$rows = @(
    @{ Row=4; A=6321;  B="Anous albivitta";     C="Grey Ternlet";      D="Least Concern"; E="Not threatened; Marine; Not migratory";                              F="Not listed"; H="1_2_3_4_5_6"; I=2005; J=1700; K="both"; M="TRUE"; N="A"; O=406 }
    @{ Row=3; A=28322; B="Aplonis metallica";   C="Metallic Starling"; D="Least Concern"; E="Not threatened; Marine; Not migratory";                              F="Not listed"; H="1_2_3_4_5_6"; I=2005; J=1700; K="land"; M="TRUE"; N="A"; O=1935 }
    @{ Row=5; A=5706;  B="Charadrius veredus";  C="Oriental Plover";   D="Least Concern"; E="Not threatened; Marine; Migratory (Bonn, CAMBA, JAMBA, ROKAMBA)";     F="Not listed"; H="1_2_3_4_5_6"; I=2005; J=1700; K="land"; M="TRUE"; N="A"; O=306 }
    @{ Row=6; A=6716;  B="Pachyptila salvini";  C="Salvin's Prion";    D="Least Concern"; E="Not threatened; Marine; Not migratory";                              F="Not listed"; H="1_2_3_4_5_6"; I=2005; J=1700; K="both"; M="TRUE"; N="A"; O=544 }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = $r.A   # A
    $ws.Cells.Item($row, 2).Value = $r.B   # B
    $ws.Cells.Item($row, 3).Value = $r.C   # C
    $ws.Cells.Item($row, 4).Value = $r.D   # D
    $ws.Cells.Item($row, 5).Value = $r.E   # E
    $ws.Cells.Item($row, 6).Value = $r.F   # F
    $ws.Cells.Item($row, 8).Value = $r.H   # H
    $ws.Cells.Item($row, 9).Value = $r.I   # I
    $ws.Cells.Item($row, 10).Value = $r.J  # J
    $ws.Cells.Item($row, 11).Value = $r.K  # K
    # Column M holds the text value "TRUE" (quoted as text, like row 1/2) -
    # copy the already-formatted cell so it keeps the shared "TRUE" string
    # and the text number-format style instead of being coerced to a boolean.
    $ws.Range("M1").Copy($ws.Cells.Item($row, 13))
    $ws.Cells.Item($row, 14).Value = $r.N  # N
    $ws.Cells.Item($row, 15).Value = $r.O  # O
}

$ws.Range("S9").Select()

Write-Output "done"
